# Add a new "RTS_CTS" worksheet at the end of the workbook with the first
# partial batch of RTS/CTS results (commit: "Added first results wit RTS/CTS").

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (and becomes the active sheet, like in the authored file).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "RTS_CTS"

# Header row: packet-size labels, same layout used by the other result sheets.
$ws.Range("C3").Value = 333
$ws.Range("D3").Value = 334
$ws.Range("E3").Value = 335
$ws.Range("F3").Value = 336
$ws.Range("G3").Value = 337

# "Number of stations" column.
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 7
$ws.Range("B6").Value = 14
$ws.Range("B7").Value = 21

# Only the run for 21 stations has data so far; the rest are still empty
# (hence the #DIV/0! from AVERAGE/STDEV over blank ranges).
$ws.Range("C7").Value = 398.99234772000011
$ws.Range("E7").Value = 349.52111785999978
$ws.Range("F7").Value = 425.85943592000001
$ws.Range("G7").Value = 428.08839937000022

# Summary formulas (average / stdev / 95% confidence interval) for every row.
$ws.Range("H4:H6").Formula = "=AVERAGE(C4:G4)"
$ws.Range("I4:I6").Formula = "=STDEV(C4:G4)"
$ws.Range("J4:J6").Formula = "=_xlfn.CONFIDENCE.NORM(0.05, I4, COUNTA(C4:G4))"

$ws.Range("H7").Formula = "=AVERAGE(C7:G7)"
$ws.Range("I7").Formula = "=STDEV(C7:G7)"
$ws.Range("J7").Formula = "=_xlfn.CONFIDENCE.NORM(0.05, I7, COUNTA(C7:G7))"

# Match the authored file's view state: RTS_CTS tab selected, K28 highlighted.
$ws.Range("K28").Select()
